$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1759259259259259
$ws.Range("C2").Value = 0.5694444444444444
$ws.Range("J2").Value = 0.03240740740740741
$ws.Range("P2").Value = 0.1296296296296296
$ws.Range("S2").Value = 0.09259259259259259
$ws.Range("C3").Value = 0.046875
$ws.Range("J3").Value = 0.015625
$ws.Range("P3").Value = 0.7578125
$ws.Range("S3").Value = 0.1796875
$ws.Range("J4").Value = 0.1162790697674419
$ws.Range("P4").Value = 0.6744186046511628
$ws.Range("S4").Value = 0.2093023255813954
$ws.Range("B6").Value = 0.06701030927835051
$ws.Range("D6").Value = 0.0154639175257732
$ws.Range("E6").Value = 0.005154639175257732
$ws.Range("F6").Value = 0.06185567010309279
$ws.Range("J6").Value = 0.3041237113402062
$ws.Range("O6").Value = 0.02577319587628866
$ws.Range("Q6").Value = 0.1391752577319588
$ws.Range("R6").Value = 0.04123711340206185
$ws.Range("S6").Value = 0.3402061855670103
$ws.Range("B7").Value = 0.09497206703910614
$ws.Range("D7").Value = 0.01675977653631285
$ws.Range("E7").Value = 0.00558659217877095
$ws.Range("F7").Value = 0.07262569832402235
$ws.Range("J7").Value = 0.1675977653631285
$ws.Range("O7").Value = 0.0223463687150838
$ws.Range("Q7").Value = 0.2569832402234637
$ws.Range("R7").Value = 0.07262569832402235
$ws.Range("S7").Value = 0.2905027932960894
$ws.Range("B8").Value = 0.08683473389355742
$ws.Range("D8").Value = 0.0196078431372549
$ws.Range("E8").Value = 0.002801120448179272
$ws.Range("F8").Value = 0.08403361344537816
$ws.Range("J8").Value = 0.0896358543417367
$ws.Range("O8").Value = 0.02240896358543417
$ws.Range("Q8").Value = 0.2100840336134454
$ws.Range("R8").Value = 0.06162464985994398
$ws.Range("S8").Value = 0.42296918767507
$ws.Range("B9").Value = 0.06542056074766354
$ws.Range("D9").Value = 0.01869158878504673
$ws.Range("F9").Value = 0.05607476635514019
$ws.Range("J9").Value = 0.102803738317757
$ws.Range("O9").Value = 0.009345794392523364
$ws.Range("Q9").Value = 0.2383177570093458
$ws.Range("R9").Value = 0.102803738317757
$ws.Range("S9").Value = 0.4065420560747663
$ws.Range("B10").Value = 0.08710801393728224
$ws.Range("D10").Value = 0.02264808362369338
$ws.Range("F10").Value = 0.06097560975609756
$ws.Range("J10").Value = 0.132404181184669
$ws.Range("O10").Value = 0.01916376306620209
$ws.Range("Q10").Value = 0.2290940766550523
$ws.Range("R10").Value = 0.07142857142857142
$ws.Range("S10").Value = 0.3771777003484321
$ws.Range("F11").Value = 0.003448275862068965
$ws.Range("G11").Value = 0.1724137931034483
$ws.Range("J11").Value = 0.09310344827586207
$ws.Range("K11").Value = 0.2517241379310345
$ws.Range("L11").Value = 0.4482758620689655
$ws.Range("S11").Value = 0.03103448275862069
$ws.Range("G12").Value = 0.7647058823529411
$ws.Range("J12").Value = 0.125
$ws.Range("K12").Value = 0.02941176470588235
$ws.Range("L12").Value = 0.03676470588235294
$ws.Range("S12").Value = 0.04411764705882353
$ws.Range("G13").Value = 0.7368421052631579
$ws.Range("J13").Value = 0.2368421052631579
$ws.Range("S13").Value = 0.02631578947368421
$ws.Range("F15").Value = 0.02617801047120419
$ws.Range("H15").Value = 0.1204188481675393
$ws.Range("I15").Value = 0.1099476439790576
$ws.Range("J15").Value = 0.3403141361256545
$ws.Range("K15").Value = 0.06282722513089005
$ws.Range("M15").Value = 0.005235602094240838
$ws.Range("O15").Value = 0.08900523560209424
$ws.Range("S15").Value = 0.2460732984293194
$ws.Range("F16").Value = 0.02666666666666667
$ws.Range("H16").Value = 0.1533333333333333
$ws.Range("I16").Value = 0.08666666666666667
$ws.Range("J16").Value = 0.4466666666666667
$ws.Range("K16").Value = 0.07333333333333333
$ws.Range("M16").Value = 0.02
$ws.Range("N16").Value = 0.006666666666666667
$ws.Range("O16").Value = 0.08
$ws.Range("S16").Value = 0.1066666666666667
$ws.Range("F17").Value = 0.02391304347826087
$ws.Range("H17").Value = 0.1304347826086956
$ws.Range("I17").Value = 0.1130434782608696
$ws.Range("J17").Value = 0.4565217391304348
$ws.Range("K17").Value = 0.09130434782608696
$ws.Range("M17").Value = 0.01304347826086956
$ws.Range("N17").Value = 0.002173913043478261
$ws.Range("O17").Value = 0.05869565217391304
$ws.Range("S17").Value = 0.1108695652173913
$ws.Range("F18").Value = 0.006802721088435374
$ws.Range("H18").Value = 0.1360544217687075
$ws.Range("I18").Value = 0.1156462585034014
$ws.Range("J18").Value = 0.4421768707482993
$ws.Range("K18").Value = 0.1564625850340136
$ws.Range("M18").Value = 0.01360544217687075
$ws.Range("O18").Value = 0.05442176870748299
$ws.Range("S18").Value = 0.07482993197278912
$ws.Range("F19").Value = 0.01576994434137291
$ws.Range("H19").Value = 0.2077922077922078
$ws.Range("I19").Value = 0.1020408163265306
$ws.Range("J19").Value = 0.3794063079777366
$ws.Range("K19").Value = 0.1094619666048238
$ws.Range("M19").Value = 0.02504638218923933
$ws.Range("N19").Value = 0.0009276437847866419
$ws.Range("O19").Value = 0.0575139146567718
$ws.Range("S19").Value = 0.1020408163265306
